$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing strings (text changes)
$ws.Range("K28").Value = "Stop avatar jumping around"
$ws.Range("L28").Value = "Change end without reposistioning avatar"

# Add new row 29
$ws.Range("L29").Value = "Change start without moving avatar"
$ws.Range("M29").Value = "DONE"
$ws.Range("M24").Copy()
$ws.Range("M29").PasteSpecial(-4122)

# Match the final selection left by the author (active cell on the new last row)
$ws.Range("M29").Select() | Out-Null

